$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing existing rows 6-11 down to 7-12
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new weekly price observation
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44428
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112003
$ws.Cells.Item(6, 7).Value = "Ajo"
$ws.Cells.Item(6, 8).Value = "Chino"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 300
$ws.Cells.Item(6, 11).Value = 15000
$ws.Cells.Item(6, 12).Value = 16000
$ws.Cells.Item(6, 13).Value = 15500
$ws.Cells.Item(6, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(6, 15).Value = "China"
$ws.Cells.Item(6, 16).Value = 1550
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = "Hortaliza"
